# Update "想去人数" (F column) values on both the "展览" and "全部类型"
# worksheets, which contain duplicated data for the same set of events.

$wb = $excel.ActiveWorkbook

# Map of worksheet name -> hashtable of row number -> new F value
$updates = @{
    "展览" = @{
        4  = 3439
        5  = 233
        6  = 4948
        7  = 493
        9  = 184
        10 = 655
        12 = 60
        16 = 29
        21 = 4813
        22 = 34
        25 = 5949
        28 = 3208
        30 = 690
        33 = 106
        35 = 920
        40 = 913
    }
    "全部类型" = @{
        8  = 3439
        9  = 233
        10 = 4948
        11 = 493
        13 = 184
        14 = 655
        16 = 60
        20 = 29
        26 = 4813
        27 = 34
        30 = 5949
        33 = 3208
        35 = 690
        39 = 106
        40 = 920
        45 = 913
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($row in $rows.Keys) {
        $ws.Range("F$row").Value = $rows[$row]
    }
}
